# Коммуникационный план — add "Управление коммуникацией с каждой группой" block
# (5 new task rows), per commit "add 5 task (fish)".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Section title row (12): merged B12:F12, centered, no border/fill
# (same look as the sheet's top title in row 1).
# ---------------------------------------------------------------------------
$ws.Range("B1").Copy()
$ws.Range("B12:F12").PasteSpecial(-4122)
$ws.Range("B12").Value = "Управление коммуникацией с каждой группой"
$ws.Range("B12:F12").Merge()
$ws.Rows.Item(12).RowHeight = 14.25

# ---------------------------------------------------------------------------
# Column header row (13): B13 blank w/ plain border, C13:F13 bordered /
# centered / wrapped with a grey header fill.
# ---------------------------------------------------------------------------
$ws.Range("C4").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("B13").HorizontalAlignment = 1
$ws.Range("B13").VerticalAlignment = -4107
$ws.Range("B13").WrapText = $false

$ws.Range("C4").Copy()
$ws.Range("C13:F13").PasteSpecial(-4122)
$ws.Range("C13:F13").Interior.Color = 13553360

$ws.Range("C13").Value = "Какая информация требуется участнику или группе?"
$ws.Range("D13").Value = "Кто будет её передавать?"
$ws.Range("E13").Value = "Как часто нужно её передавать?"
$ws.Range("F13").Value = "По каким каналам связи её будут передавать?"
$ws.Rows.Item(13).RowHeight = 42.75

# ---------------------------------------------------------------------------
# Data rows 14-18, one per stakeholder group. Formatting is copied from the
# matching stakeholder rows above (4-8) so the look (borders / fills /
# wrap / centering) matches exactly; only the text differs.
# ---------------------------------------------------------------------------

# Row 14 — Project sponsor
$ws.Range("B4").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("D4").Copy()
$ws.Range("C14:D14").PasteSpecial(-4122)
$ws.Range("C4").Copy()
$ws.Range("E14:F14").PasteSpecial(-4122)
$ws.Range("B14").Value = "Project sponsor"
$ws.Range("C14").Value = "Статус выполнения проекта"
$ws.Range("D14").Value = "Product manager"
$ws.Range("E14").Value = "по Agile - каждые 2 недели"
$ws.Range("F14").Value = "Личная встреча"
$ws.Rows.Item(14).RowHeight = 14.25

# Row 15 — Product manager
$ws.Range("B5").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("E5").Copy()
$ws.Range("C15:F15").PasteSpecial(-4122)
$ws.Range("B15").Value = "Product manager"
$ws.Range("C15").Value = "Какие ранне поставленные задачи выполнены"
$ws.Range("D15").Value = "Система учета (WEEEK.ru) и митинги"
$ws.Range("E15").Value = "Раз в день для корректировки планов при возникновеннии проблем"
$ws.Range("F15").Value = "Онлайн и offline встречи"
$ws.Rows.Item(15).RowHeight = 42.75

# Row 16 — Маркетолог
$ws.Range("B6").Copy()
$ws.Range("B16").PasteSpecial(-4122)
$ws.Range("C6").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("D6").Copy()
$ws.Range("D16:E16").PasteSpecial(-4122)
$ws.Range("C6").Copy()
$ws.Range("F16").PasteSpecial(-4122)
$ws.Range("B16").Value = "Маркетолог"
$ws.Range("C16").Value = "Описание проекта, цель проекта, реализованные функции, графические материалы"
$ws.Range("D16").Value = "Product manager"
$ws.Range("E16").Value = "по Agile - каждые 2 недели"
$ws.Range("F16").Value = "Через мессенджер"
$ws.Rows.Item(16).RowHeight = 57

# Row 17 — Developers
$ws.Range("B7").Copy()
$ws.Range("B17").PasteSpecial(-4122)
$ws.Range("C7").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("D7").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("C7").Copy()
$ws.Range("E17").PasteSpecial(-4122)
$ws.Range("D7").Copy()
$ws.Range("F17").PasteSpecial(-4122)
$ws.Range("B17").Value = "Developers"
$ws.Range("C17").Value = "Списко задач с временем отведенным на их выполнение"
$ws.Range("D17").Value = "Product manager"
$ws.Range("E17").Value = "Раз в 2 недели, но необходима промежуточная корректировка"
$ws.Range("F17").Value = "Система учета (WEEEK.ru) и митинги"
$ws.Rows.Item(17).RowHeight = 42.75

# Row 18 — Пользователи
$ws.Range("B8").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("E8").Copy()
$ws.Range("C18:E18").PasteSpecial(-4122)
$ws.Range("F8").Copy()
$ws.Range("F18").PasteSpecial(-4122)
$ws.Range("B18").Value = "Пользователи"
$ws.Range("C18").Value = "Описание проекта, новые функции и преимущества"
$ws.Range("D18").Value = "Макетолог"
$ws.Range("E18").Value = "Раз в 2 недели, но можно чаще, при исправлении багов или оповешении пользователей о выгодных предложениях"
$ws.Range("F18").Value = "Социальные сети, реклама на других сайтах"
$ws.Rows.Item(18).RowHeight = 57

Write-Output "edit complete"
